$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description: ..." paragraph right after the Heading1
#    title paragraph ("Play 100 Zombies Free: Review, Pros and Cons").
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs.Item(2)
$metaXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' +
  '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' +
  '<w:r><w:t>: Explore the features of 100 Zombies and decide if it is worth playing for free. Read our review covering gameplay, design, symbols, and more.</w:t></w:r></w:p>' +
  '</w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'
$metaPara.Range.InsertXML($metaXml)

# ---------------------------------------------------------------------------
# 2) Remove the trailing duplicate title paragraph ("Play 100 Zombies Free:
#    Review, Pros and Cons") near the end of the document, and replace the
#    text of the following (italic) paragraph with the new image prompt.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs.Item($count - 1)
$descPara = $d.Paragraphs.Item($count)

$deleteRange = $d.Range($dupTitlePara.Range.Start, $descPara.Range.Start)
$deleteRange.Delete()

$count2 = $d.Paragraphs.Count
$descPara2 = $d.Paragraphs.Item($count2)
$descRange = $d.Range($descPara2.Range.Start, $descPara2.Range.End)
$descRange.Text = 'Create a feature image for "100 Zombies" game that features a happy Maya warrior with glasses in cartoon style. Maya warrior should be depicted wearing protective gear and holding a rifle as if ready to face a zombie outbreak. The background should resemble a city engulfed in flames with expanding clouds to emphasize the unhealthy environment caused by the epidemic. The overall style of the image should be fun and playful, capturing the game''s unique take on the popular zombie theme.'
